# feat: add video to dictionary entry format
#
# Inserts two new columns (video_desc / video_fn) immediately before the
# existing theme / secondary_theme columns on both "Sheet1" and "Skip",
# fills in the new values for the rows that already have audio (the first
# three dictionary entries), updates the named ranges that describe the
# data tables, and restores selection/active-sheet state.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Skip")

# --- Sheet1 ------------------------------------------------------------
# Insert two blank columns at F:G, pushing the existing theme (F) /
# secondary_theme (G) columns to H:I.
$ws1.Columns("F:G").Insert()

$ws1.Range("F1").Value = "Nolan Van Hell"
$ws1.Range("G1").Value = "snowfall.mp4"

$ws1.Range("F2").Value = "Nolan Van Hell"
$ws1.Range("G2").Value = "snowfall.mp4"

$ws1.Range("F3").Value = "Nolan Van Hell"
$ws1.Range("G3").Value = "snowfall.mp4"

# Row 4 ("goodbye") never had audio/video columns filled in, so F4/G4 stay
# blank, same as D4/E4 already were.

$ws1.Columns("F").ColumnWidth = 12.5
$ws1.Columns("G").ColumnWidth = 11.5

# --- Sheet2 (Skip) -------------------------------------------------------
$ws2.Columns("F:G").Insert()

$ws2.Range("F1").Value = "video_desc"
$ws2.Range("G1").Value = "video_fn"

$ws2.Range("F2").Value = "Nolan Van Hell"
$ws2.Range("G2").Value = "snowfall.mp4"

$ws2.Range("F3").Value = "Nolan Van Hell"
$ws2.Range("G3").Value = "snowfall.mp4"

$ws2.Range("F4").Value = "Nolan Van Hell"
$ws2.Range("G4").Value = "snowfall.mp4"

# Row 5 ("goodbye") stays without video columns, mirroring Sheet1 row 4.

$ws2.Columns("F").ColumnWidth = 12.5
$ws2.Columns("G").ColumnWidth = 11.5

# --- Defined names -------------------------------------------------------
# The "data" / "data_1" named ranges describe the full table extents; they
# need to grow from column G to column I to cover the two new columns.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "Sheet1!data") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$I`$4"
    }
    if ($n.Name -eq "Skip!data_1") {
        $n.RefersTo = "=Skip!`$A`$2:`$I`$5"
    }
}

# --- Selection / active sheet --------------------------------------------
# Restore per-sheet selection; select Sheet2 first so that Sheet1 ends up
# as the active (tab-selected) sheet, matching the original workbook.
[void]$ws2.Range("H10").Select()
[void]$ws1.Range("F12").Select()
